# Refactor quit confirmation logic for improved user experience and add
# combat handling in the App class — bug tracker sheet gains a new row
# documenting the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bug-report row (row 11)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "quit command accepts characters other than y or n"
$ws.Range("C11").Value = "ricky"
$ws.Range("D11").Value = "confirmation doesn't enforce y or n"
$ws.Range("E11").Value = "Wrap the confirmation in a loop that continues prompting until y or n"
$ws.Range("F11").Value = "fixed"

# Column E now holds longer text, so it gets re-sized to fit the content.
$ws.Columns.Item(5).ColumnWidth = 62

# The active selection ends up one row below/right of the new data.
$ws.Range("F12").Select()
